# This workbook contains weekly price records for Fruta / Femacal de La
# Calera - Papaya. The edit reshuffles the per-record fields (date,
# quality, volume, min/max/avg price, trade unit, $/Kg price and Kg per
# unit) across the existing data rows, while leaving the row's fixed
# market/product metadata (columns A, B, C, E, F, G, H, I, J, K, R)
# untouched.
#
# Mapping below: destination row -> source row (both refer to row numbers
# in the ORIGINAL worksheet, i.e. before any writes happen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 60
    3  = 56
    4  = 61
    5  = 43
    6  = 44
    7  = 62
    8  = 63
    9  = 64
    10 = 21
    11 = 22
    12 = 55
    13 = 23
    14 = 57
    15 = 28
    16 = 29
    17 = 7
    18 = 15
    19 = 16
    20 = 54
    21 = 67
    22 = 9
    23 = 17
    24 = 18
    25 = 41
    26 = 42
    27 = 37
    28 = 35
    29 = 36
    30 = 19
    31 = 5
    32 = 6
    33 = 53
    34 = 48
    35 = 31
    36 = 32
    37 = 30
    38 = 38
    39 = 49
    40 = 50
    41 = 12
    42 = 2
    43 = 66
    44 = 14
    45 = 47
    46 = 40
    47 = 39
    48 = 51
    49 = 26
    50 = 45
    51 = 46
    52 = 25
    53 = 13
    54 = 52
    55 = 27
    56 = 10
    57 = 11
    58 = 65
    59 = 3
    60 = 4
    61 = 20
    62 = 8
    63 = 34
    64 = 24
    65 = 58
    66 = 59
    67 = 33
}

# Columns whose values get permuted per row.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot the original value of every affected cell before writing
# anything, since source and destination rows overlap.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 67; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot["$col$srcRow"]
    }
}
